# Add a new worksheet "Лист2" (populated with its own data) after the
# existing sheet, and make it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Copy the first sheet to the end of the workbook, so the new sheet starts
# out with the same formatting context (row height, etc.) as the rest of
# the workbook, then strip its copied content and rebuild it from scratch.
$sheet1 = $wb.Worksheets.Item(1)
$sheet1.Copy($null, $sheet1)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Лист2"
$newSheet.Cells.Clear()

# Populate data on the new sheet.
$newSheet.Range("A1").Value = "l21"
$newSheet.Range("B3").Value = "l22"
$newSheet.Range("C4").Value = "l23"
$newSheet.Range("D5").Value = "l24"
$newSheet.Range("E5").Value = "l25"

# Select a cell on the new sheet to mimic the authored selection state.
$newSheet.Range("E6").Select()

# Make the new sheet the active tab.
$newSheet.Activate()
